$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Swap the values of P8 and Q8
$p8 = $ws.Range("P8").Value()
$q8 = $ws.Range("Q8").Value()
$ws.Range("P8").Value = $q8
$ws.Range("Q8").Value = $p8

# Update the active cell / selection to P8
$ws.Activate()
$ws.Range("P8").Select()
